$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$lo = $ws.ListObjects.Item(1)
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("C1:C31"), 0, 2)
$lo.Sort.Header = 1
$lo.Sort.Apply()
$ws.Range("B8").Select()
